# Optuna Attempt (go back with original)
# Update the "Forecast Comparison" sheet with revised forecast/metric values
# and the "Summary" sheet with the resulting aggregate figures.

$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison sheet updates ---

# Row 2
$wsForecast.Range("D2").Value = 52
$wsForecast.Range("H2").Value = 7.26
$wsForecast.Range("L2").Value = 1.12

# Row 3
$wsForecast.Range("D3").Value = 56
$wsForecast.Range("H3").Value = 5.87
$wsForecast.Range("L3").Value = 1.19

# Row 4
$wsForecast.Range("D4").Value = 47
$wsForecast.Range("H4").Value = 5.72

# Row 5
$wsForecast.Range("D5").Value = 47
$wsForecast.Range("H5").Value = 4.72
$wsForecast.Range("L5").Value = 1.08

# Row 6
$wsForecast.Range("D6").Value = 49
$wsForecast.Range("H6").Value = 3.62
$wsForecast.Range("I6").Value = "Low"
$wsForecast.Range("J6").Value = "Normal"
$wsForecast.Range("L6").Value = 0.83

# Row 7
$wsForecast.Range("H7").Value = 1.52
$wsForecast.Range("I7").Value = "Low"
$wsForecast.Range("J7").Value = "Normal"
$wsForecast.Range("L7").Value = 1

# Row 8
$wsForecast.Range("H8").Value = 0.46
$wsForecast.Range("L8").Value = 0.91

# Row 9
$wsForecast.Range("D9").Value = 52
$wsForecast.Range("L9").Value = 1.07

# Row 10
$wsForecast.Range("D10").Value = 46
$wsForecast.Range("L10").Value = 1

# Row 11
$wsForecast.Range("D11").Value = 49
$wsForecast.Range("L11").Value = 1.01

# Row 12
$wsForecast.Range("D12").Value = 49
$wsForecast.Range("L12").Value = 0.87

# Row 13
$wsForecast.Range("L13").Value = 1

# Row 14
$wsForecast.Range("L14").Value = 0.99

# Row 15
$wsForecast.Range("D15").Value = 46
$wsForecast.Range("L15").Value = 1.14

# Row 16
$wsForecast.Range("D16").Value = 49
$wsForecast.Range("L16").Value = 0.82

# Row 17
$wsForecast.Range("D17").Value = 46
$wsForecast.Range("L17").Value = 0.9

# --- Summary sheet updates ---

$wsSummary.Range("B9").Value  = "961"
$wsSummary.Range("B10").Value = "484"
$wsSummary.Range("B11").Value = "204"
$wsSummary.Range("B12").Value = "95"
$wsSummary.Range("B14").Value = "47"
